# Actualizacion de precios junio 2025
# cotizador.xlsx - Hoja1: +2% price update across the whole price grid
# (B2:H11) plus the incidental formatting/selection state Excel leaves
# behind when the sheet is re-saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-touch the header font (same name, same value) - this is what nudges
# Excel into (re)writing the <family val="2"/> element for that font the
# next time the workbook is saved.
$ws.Range("A1:H1").Font.Name = "Calibri"

# Re-affirm protection/number-format on the price cells so the cellXf
# picks up applyNumberFormat/applyProtection (same currency format as
# before - only the bookkeeping flags change on resave).
$priceRange = $ws.Range("B2:H11")
$priceRange.NumberFormat = '_ "$"\ * #,##0.00_ ;_ "$"\ * \-#,##0.00_ ;_ "$"\ * "-"??_ ;_ @_ '
$priceRange.Locked = $true

# --- Price updates (+2%) ---------------------------------------------
# Row 2 (18-25)
$ws.Range("B2").Value = 94108
$ws.Range("C2").Value = 134720
$ws.Range("D2").Value = 159074
$ws.Range("E2").Value = 269677
$ws.Range("F2").Value = 418979
$ws.Range("G2").Value = 81963
$ws.Range("H2").Value = 70699

# Row 3 (26-35)
$ws.Range("B3").Value = 94519
$ws.Range("C3").Value = 136792
$ws.Range("D3").Value = 160580
$ws.Range("E3").Value = 277101
$ws.Range("F3").Value = 420277
$ws.Range("G3").Value = 82177
$ws.Range("H3").Value = 70978

# Row 4 (36-54)
$ws.Range("B4").Value = 106324
$ws.Range("C4").Value = 154206
$ws.Range("D4").Value = 180943
$ws.Range("E4").Value = 311894
$ws.Range("F4").Value = 473215
$ws.Range("G4").Value = 92349
$ws.Range("H4").Value = 80060

# Row 5 (55-59)
$ws.Range("B5").Value = 135071
$ws.Range("C5").Value = 193128
$ws.Range("D5").Value = 226212
$ws.Range("E5").Value = 379076
$ws.Range("F5").Value = 564445
$ws.Range("G5").Value = 117828
$ws.Range("H5").Value = 102488

# Row 6
$ws.Range("B6").Value = 232249
$ws.Range("C6").Value = 309367
$ws.Range("D6").Value = 407177
$ws.Range("E6").Value = 567532
$ws.Range("F6").Value = 773392
$ws.Range("G6").Value = 203503
$ws.Range("H6").Value = 177608

# Row 7
$ws.Range("B7").Value = 60738
$ws.Range("C7").Value = 87540
$ws.Range("D7").Value = 108697
$ws.Range("E7").Value = 205447
$ws.Range("F7").Value = 299739
$ws.Range("G7").Value = 57737
$ws.Range("H7").Value = 51963

# Row 8
$ws.Range("B8").Value = 61878
$ws.Range("C8").Value = 88404
$ws.Range("D8").Value = 108661
$ws.Range("E8").Value = 207139
$ws.Range("F8").Value = 305303
$ws.Range("G8").Value = 57461
$ws.Range("H8").Value = 53631

# Row 9
$ws.Range("B9").Value = 71019
$ws.Range("C9").Value = 100072
$ws.Range("D9").Value = 125458
$ws.Range("E9").Value = 240333
$ws.Range("F9").Value = 352521
$ws.Range("G9").Value = 67179
$ws.Range("H9").Value = 62927

# Row 10
$ws.Range("B10").Value = 86553
$ws.Range("C10").Value = 121672
$ws.Range("D10").Value = 149432
$ws.Range("E10").Value = 289346
$ws.Range("F10").Value = 412508
$ws.Range("G10").Value = 80706
$ws.Range("H10").Value = 75598

# Row 11
$ws.Range("B11").Value = 79991
$ws.Range("C11").Value = 114511
$ws.Range("D11").Value = 135212
$ws.Range("E11").Value = 229226
$ws.Range("F11").Value = 356131
$ws.Range("G11").Value = 69669
$ws.Range("H11").Value = 60094

# Leave the selection where the author's last save left it.
$ws.Range("D14").Select() | Out-Null
